$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (A=6): name shifts to "line7"; new random C/D; now in service
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# Row 9 (A=7): name shifts to "line8"; new random C/D; now in service
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Row 10 (A=8): name shifts to "extr1"; C/D take prior row-8 values
$ws.Range("B10").Value = "extr1"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12

# Row 11 (A=9): name shifts to "extr2"; C/D take prior row-9 values
$ws.Range("B11").Value = "extr2"
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9

# Row 12 (A=10): name shifts to "extr3"; C takes prior row-10 value (D unchanged)
$ws.Range("B12").Value = "extr3"
$ws.Range("C12").Value = 10

# Row 13 (A=11): name shifts to "extr4"; D takes prior row-11 value (C unchanged)
$ws.Range("B13").Value = "extr4"
$ws.Range("D13").Value = 8

# Row 14 (A=12): name shifts to "extr5"; C/D take prior row-12 values
$ws.Range("B14").Value = "extr5"
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11

# Row 15 (A=13): name shifts to "extr6"; C/D take prior row-13 values
$ws.Range("B15").Value = "extr6"
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11

# New rows 16 and 17: copy formatting from row 15 (A column bold/border style),
# then fill in values.
$ws.Range("A15:E15").Copy()
$ws.Range("A16:E17").PasteSpecial(-4122)

# New row 16 (A=14): name "extr7"; C/D take prior row-14 values
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $false

# New row 17 (A=15): name "extr8"; C/D take prior row-15 values
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false
